# Applies the "novo lm" (new LM) retraining results:
#  - column A model labels are reshuffled across rows 2..26
#  - columns B..Q (the metric values) are overwritten with the
#    freshly computed metrics, identical for every data row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new model name in column A
$rowNames = @{
    2 = "model_20_8_0"
    3 = "model_20_8_22"
    4 = "model_20_8_21"
    5 = "model_20_8_20"
    6 = "model_20_8_19"
    7 = "model_20_8_18"
    8 = "model_20_8_17"
    9 = "model_20_8_16"
    10 = "model_20_8_15"
    11 = "model_20_8_14"
    12 = "model_20_8_13"
    13 = "model_20_8_23"
    14 = "model_20_8_12"
    15 = "model_20_8_10"
    16 = "model_20_8_9"
    17 = "model_20_8_8"
    18 = "model_20_8_7"
    19 = "model_20_8_6"
    20 = "model_20_8_5"
    21 = "model_20_8_4"
    22 = "model_20_8_3"
    23 = "model_20_8_2"
    24 = "model_20_8_1"
    25 = "model_20_8_11"
    26 = "model_20_8_24"
}

# New metric values (columns B:Q), shared by every data row
$colValues = @{
    "B" = 0.9999805300082542
    "C" = 0.9991182316315311
    "D" = 0.9999729332048588
    "E" = 0.9999882717448302
    "F" = 0.9999787623451173
    "G" = 0.00001817438266093276
    "H" = 0.0008230920668076507
    "I" = 0.00004471350215755075
    "J" = 0.000001989184062598218
    "K" = 0.00002335134311007448
    "L" = 0.0002731515751168351
    "M" = 0.004263142345844525
    "N" = 1.000035944600146
    "O" = 0.004444633141428372
    "P" = 95.83099500301532
    "Q" = 140.9294005231387
}

foreach ($row in $rowNames.Keys) {
    $ws.Range("A" + $row).Value = $rowNames[$row]
    foreach ($col in $colValues.Keys) {
        $ws.Range($col + $row).Value = $colValues[$col]
    }
}

"edit complete"
